# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.946.11"
$ws.Range("E2").Value = "  +7.93%  "
$ws.Range("D3").Value = "2.425.23"
$ws.Range("E3").Value = "  +5.95%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.Formula = "'114.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +11.91%  "
$c = $ws.Range("D6")
$c.Formula = "'319.69"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.61%  "
$c = $ws.Range("D7")
$c.Formula = "'0.636"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.55%  "
$ws.Range("E8").Value = "  +0.04%  "
$c = $ws.Range("D9")
$c.Formula = "'0.634"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.26%  "
$c = $ws.Range("D10")
$c.Formula = "'43.20"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +10.85%  "
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("E12").Value = "  +6.83%  "
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("E14").Value = "  +2.33%  "
$c = $ws.Range("D15")
$c.Formula = "'16.03"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.35%  "
$ws.Range("D16").Value = "2.791.83"
$ws.Range("E16").Value = "  +5.97%  "
$ws.Range("D17").Value = "2.421.61"
$ws.Range("E17").Value = "  +5.74%  "
$ws.Range("D18").Value = "45.889.71"
$ws.Range("E18").Value = "  +7.93%  "
$c = $ws.Range("D19")
$c.Formula = "'7.66"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("E20").Value = "  +4.94%  "
$c = $ws.Range("D21")
$c.Formula = "'13.47"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "
$c = $ws.Range("D22")
$c.Formula = "'75.46"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.88%  "
$c = $ws.Range("D23")
$c.Formula = "'3.55"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +5.10%  "
$c = $ws.Range("D24")
$c.Formula = "'269.51"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  +8.53%  "
$ws.Range("E26").Value = "  -0.51%  "
$c = $ws.Range("D27")
$c.Formula = "'7.73"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +7.43%  "
$c = $ws.Range("D28")
$c.Formula = "'11.41"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.90%  "
$ws.Range("E29").Value = "  +3.77%  "
$c = $ws.Range("D30")
$c.Formula = "'39.41"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +10.60%  "
$ws.Range("E31").Value = "  +3.24%  "
$c = $ws.Range("D32")
$c.Formula = "'0.0974"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +14.40%  "
$c = $ws.Range("D33")
$c.Formula = "'173.62"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("E34").Value = "  +16.90%  "
$ws.Range("E35").Value = "  +9.76%  "
$ws.Range("E36").Value = "  +1.93%  "
$c = $ws.Range("D37")
$c.Formula = "'4.99"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +10.35%  "
$c = $ws.Range("D38")
$c.Formula = "'4.18"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +15.77%  "
$c = $ws.Range("D39")
$c.Formula = "'3.13"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +12.28%  "
$ws.Range("E40").Value = "  +6.07%  "
$ws.Range("E41").Value = "  +16.01%  "
$c = $ws.Range("D42")
$c.Formula = "'102.58"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("E43").Value = "  +6.34%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Range("D44")
$c.Formula = "'72.50"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D45")
$c.Formula = "'13.55"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +12.15%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D47")
$c.Formula = "'5.87"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +14.10%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D48")
$c.Formula = "'118.11"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +7.43%  "
$c = $ws.Range("D49")
$c.Formula = "'1.68"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Formula = "'9.46"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +9.46%  "
$c = $ws.Range("D51")
$c.Formula = "'79.24"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.36%  "
